$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-10"

# Update the "through" date header label (column I, year-to-date column)
$ws.Range("I1").Value = "2022 (through 11-10)"

# Update November 2022 carjacking count (I12) and recompute Total (I14)
$ws.Range("I12").Value = 27
$ws.Range("I14").Value = 1426
